# Add "sword world replays" entries to the checklist
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23: Rhapsody of Rogues (Sword World Replay 1)
$ws.Range("A23").Value = 1989
$ws.Range("B23").Value = "盗賊たちの狂詩曲ソード・ワールドRPGリプレイ集〈1〉"
$ws.Range("C23").Value = "Rhapsody of Rogues: Sword World Replay 1"
$ws.Range("D23").Value = "Fujimi Shobo"
$ws.Range("E23").Value = "rhapsody_of_rogues.jpg"

# Row 24: Symphony of Monsters (Sword World Replay 2)
$ws.Range("A24").Value = 1990
$ws.Range("E24").Value = "symphony_of_monsters.jpg"
$ws.Range("C24").Value = "Symphony of Monsters: Sword World Replay 2"
$ws.Range("B24").Value = "モンスターたちの交響曲ソード・ワールドRPGリプレイ集〈2〉"
$ws.Range("D24").Value = "Fujimi Shobo"

# Row 25: Endless Improvisation (Sword World Replay 3)
$ws.Range("A25").Value = 1991
$ws.Range("B25").Value = "終わりなき即興曲ソード・ワールドRPGリプレイ集〈3〉"
$ws.Range("C25").Value = "Endless Improvisation: Sword World Replay 3"
$ws.Range("D25").Value = "Fujimi Shobo"
$ws.Range("E25").Value = "endless_improvisation.jpg"

# Widen column B to fit the longer Japanese titles (XML width 60.5)
$ws.Columns.Item(2).ColumnWidth = 59.665

# Move the active selection to just past the new last row, as in the original edit
$ws.Range("A26").Select()
